# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 25 in the "Achicoria" sheet,
# pushing the existing rows 25-39 down to 26-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 25 (shifts 25..39 -> 26..40).
$ws.Rows(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value = 10
$ws.Cells.Item(25, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(25, 3).Value = "La Araucanía"
$ws.Cells.Item(25, 4).Value = 44726
$ws.Cells.Item(25, 5).Value = 9
$ws.Cells.Item(25, 6).Value = 100112010
$ws.Cells.Item(25, 7).Value = "Achicoria"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 125
$ws.Cells.Item(25, 11).Value = 8000
$ws.Cells.Item(25, 12).Value = 8000
$ws.Cells.Item(25, 13).Value = 8000
$ws.Cells.Item(25, 14).Value = "`$/caja 18 unidades"
$ws.Cells.Item(25, 15).Value = "Región Metropolitana"
$ws.Cells.Item(25, 16).Value = 444
$ws.Cells.Item(25, 17).Value = 18
$ws.Cells.Item(25, 18).Value = "Hortaliza"
